# Correccion tablero: "Se aprobo el documentacion de Arquitectura" ->
# "Se aprobo el documento de Arquitectura" (wording fix, 3rd slide,
# "TextBox 117" status box, the "Se aprobo ... Arquitectura" bullet).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item("TextBox 117")
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Locate the bullet paragraph that still has the old wording.
$targetPara = $null
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    if ($candidate.Text -like "*documentación de Arquitectura*") {
        $targetPara = $candidate
        break
    }
}

# Replace just the "el documentación " span with "el documento " so the
# run splits into: "Se aprobó " | "el documento " | "de Arquitectura ( Fecha fin: 16-09-2016)."
$fullText = $targetPara.Text
$oldSpan  = "el documentación "
$newSpan  = "el documento "
$spanStart = $fullText.IndexOf($oldSpan)

$spanRange = $targetPara.Characters($spanStart + 1, $oldSpan.Length)
$spanRange.Text = $newSpan
